# Deep Audit: Remove all hardcoded defaults from the Finance Dashboard.
# Zeroes out the hardcoded placeholder assumptions that were baked into the
# LIQUIDITY_MONITOR, PROFIT_CONTROL, and DEBT_MANAGER sheets so the model
# relies on actual inputs instead of stale sample defaults.

$wb = $excel.ActiveWorkbook

# --- LIQUIDITY_MONITOR: Procurement Spend / Receivables / Payables (Est.) ---
$wsLiquidity = $wb.Worksheets.Item("LIQUIDITY_MONITOR")
$wsLiquidity.Range("B17:I17").Value = 0   # Procurement Spend (Est.)
$wsLiquidity.Range("B19:I19").Value = 0   # Receivables (Hard)
$wsLiquidity.Range("B20:I20").Value = 0   # Payables (Hard)

# --- PROFIT_CONTROL: Interest Expense hardcoded default ---
$wsProfit = $wb.Worksheets.Item("PROFIT_CONTROL")
$wsProfit.Range("B16:C16").Value = 0      # Interest Expense

# --- DEBT_MANAGER: Loan interest-rate hardcoded defaults ---
$wsDebt = $wb.Worksheets.Item("DEBT_MANAGER")
$wsDebt.Range("C6").Value = 0             # Loan 1 Interest Rate
$wsDebt.Range("C7").Value = 0             # Loan 2 Interest Rate
$wsDebt.Range("C8").Value = 0             # Loan 3 Interest Rate
